# The workbook contains a single data table (Hortaliza, Vega Modelo de Temuco - Papa).
# A new observation row needs to be inserted as row 996, pushing the existing
# rows 996-1061 down by one (to 997-1062) while keeping all of their data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 996; this shifts rows 996:1061 down to 997:1062
# and copies the row's formatting (incl. the date number format in column D)
# from the row above, same as Excel's native "Insert Row" behaviour.
$ws.Rows.Item(996).Insert()

# Populate the newly inserted row 996 with the new record's data.
$ws.Cells.Item(996, 1).Value = 10
$ws.Cells.Item(996, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(996, 3).Value = "La Araucanía"
$ws.Cells.Item(996, 4).Value = 45013
$ws.Cells.Item(996, 5).Value = 9
$ws.Cells.Item(996, 6).Value = 100114001
$ws.Cells.Item(996, 7).Value = "Papa"
$ws.Cells.Item(996, 8).Value = "Patagonia"
$ws.Cells.Item(996, 9).Value = "1a (cosecha)"
$ws.Cells.Item(996, 10).Value = 400
$ws.Cells.Item(996, 11).Value = 12000
$ws.Cells.Item(996, 12).Value = 12000
$ws.Cells.Item(996, 13).Value = 12000
$ws.Cells.Item(996, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(996, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(996, 16).Value = 480
$ws.Cells.Item(996, 17).Value = 25
$ws.Cells.Item(996, 18).Value = "Hortaliza"
